# Update the first-name values on "test_place_order" so they carry the
# "Beta" suffix, matching the refreshed beta test data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_place_order")
$ws.Range("C2").Value = "Fn365Beta"
$ws.Range("C3").Value = "Fn513Beta"
